$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "25.403.27"
$ws.Range("E2").Value = "  -2.34%  "
$ws.Range("D3").Value = "1.573.22"
$ws.Range("E3").Value = "  -3.90%  "
$ws.Range("E4").Value = "  +0.33%  "
Set-TextValue "D5" "207.14"
$ws.Range("E5").Value = "  -3.48%  "
$ws.Range("E6").Value = "  +0.21%  "
Set-TextValue "D7" "0.479"
$ws.Range("E7").Value = "  -4.86%  "
Set-TextValue "D8" "0.245"
$ws.Range("E8").Value = "  -2.47%  "
Set-TextValue "D9" "0.0607"
$ws.Range("E9").Value = "  -2.19%  "
Set-TextValue "D10" "17.69"
$ws.Range("E10").Value = "  -3.35%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").Value = "1.804.66"
$ws.Range("E12").Value = "  -3.19%  "
$ws.Range("D13").Value = "1.585.91"
$ws.Range("E13").Value = "  -3.48%  "
Set-TextValue "D14" "4.03"
$ws.Range("E14").Value = "  -3.71%  "
Set-TextValue "D15" "0.504"
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("D16").Value = "25.412.56"
$ws.Range("E16").Value = "  -2.24%  "
Set-TextValue "D17" "59.99"
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("D18").Value = "0.0₃0707"
$ws.Range("E18").Value = "  -4.59%  "
$ws.Range("E19").Value = "  +0.31%  "
Set-TextValue "D20" "186.37"
$ws.Range("E20").Value = "  -2.47%  "
Set-TextValue "D21" "4.14"
$ws.Range("E21").Value = "  -2.37%  "
Set-TextValue "D22" "9.30"
$ws.Range("E22").Value = "  -4.03%  "
Set-TextValue "D23" "5.89"
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D25" "141.00"
$ws.Range("E25").Value = "  -2.02%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D26" "0.128"
$ws.Range("E26").Value = "  -3.62%  "
Set-TextValue "D27" "1.71"
$ws.Range("E27").Value = "  -4.39%  "
Set-TextValue "D28" "14.90"
$ws.Range("E28").Value = "  -1.86%  "
Set-TextValue "D29" "6.47"
$ws.Range("E29").Value = "  -4.54%  "
Set-TextValue "D30" "1.17"
$ws.Range("E30").Value = "  -5.47%  "
Set-TextValue "D31" "0.0463"
$ws.Range("E31").Value = "  -3.75%  "
Set-TextValue "D32" "3.05"
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("E33").Value = "  -4.26%  "
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("D36").Value = "1.083.84"
$ws.Range("E36").Value = "  -4.56%  "
$ws.Range("E37").Value = "  +0.05%  "
Set-TextValue "D38" "2.33"
$ws.Range("E38").Value = "  -4.63%  "
$ws.Range("E39").Value = "  -3.15%  "
Set-TextValue "D40" "0.773"
$ws.Range("E40").Value = "  -10.40%  "
Set-TextValue "D41" "0.492"
$ws.Range("E41").Value = "  -5.19%  "
Set-TextValue "D42" "94.67"
$ws.Range("E42").Value = "  -3.67%  "
$ws.Range("D43").Value = "1.718.14"
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("E44").Value = "  -3.04%  "
Set-TextValue "D45" "0.731"
$ws.Range("E45").Value = "  -5.65%  "
$ws.Range("D46").Value = "0.0₆0107"
$ws.Range("E46").Value = "  -7.33%  "
Set-TextValue "D47" "52.81"
$ws.Range("E47").Value = "  -3.90%  "
Set-TextValue "D48" "0.0508"
$ws.Range("E48").Value = "  -3.73%  "
$ws.Range("E49").Value = "  -4.51%  "
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("E51").Value = "  -0.04%  "
